$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 318
$ws1.Range("F4").Value = 8339
$ws1.Range("F5").Value = 6082
$ws1.Range("F6").Value = 522
$ws1.Range("F7").Value = 107
$ws1.Range("F9").Value = 71
$ws1.Range("F10").Value = 317
$ws1.Range("F11").Value = 1037
$ws1.Range("F12").Value = 81

# Sheet "演出" (Performance) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 96
$ws2.Range("F3").Value = 1
$ws2.Range("F5").Value = 7

# Sheet "全部类型" (All types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 318
$ws4.Range("F3").Value = 20
$ws4.Range("F4").Value = 8339
$ws4.Range("F5").Value = 6082
$ws4.Range("F6").Value = 522
$ws4.Range("F7").Value = 107
$ws4.Range("F8").Value = 19
$ws4.Range("F9").Value = 71
$ws4.Range("F10").Value = 317
$ws4.Range("F11").Value = 96
$ws4.Range("F12").Value = 1
$ws4.Range("F14").Value = 7
$ws4.Range("F15").Value = 1037
$ws4.Range("F16").Value = 81
$ws4.Range("F17").Value = 2

$wb.Save()
